$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from G3 onto G4 before writing the value,
# so the new row's date cell keeps the same number format (style index 1).
$ws.Range("G3").Copy()
$ws.Range("G4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A4").Value = 10035.5
$ws.Range("B4").Value = 9945
$ws.Range("C4").Value = 107.96
$ws.Range("D4").Value = 108.94
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.91
$ws.Range("G4").Value = 42609.505578703705
$ws.Range("H4").Value = $true
